# Auto-generated edit script: updates Leve profit-calculation cells
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets, matching a
# scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1003.7368
$ws.Range("J17").Value = 1043.3889
$ws.Range("L17").Value = 3130.1667
$ws.Range("N17").Value = -3466.1667
$ws.Range("H129").Value = 2091.75
$ws.Range("J129").Value = 2824.75
$ws.Range("L129").Value = 8474.25
$ws.Range("N129").Value = -18474.25
$ws.Range("H138").Value = 2197.3103
$ws.Range("I138").Value = 921.72
$ws.Range("K138").Value = 2765.16
$ws.Range("M138").Value = 2374.84

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 42577.89
$ws.Range("I32").Value = 45768.12
$ws.Range("K32").Value = 45768.12
$ws.Range("M32").Value = -45481.12
$ws.Range("H61").Value = 6927.3335
$ws.Range("I61").Value = 4258.857
$ws.Range("K61").Value = 4258.857
$ws.Range("M61").Value = -4046.857
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H74").Value = 505869.9
$ws.Range("I74").Value = 772184.7
$ws.Range("K74").Value = 772184.7
$ws.Range("M74").Value = -771310.7
$ws.Range("H77").Value = 505869.9
$ws.Range("I77").Value = 772184.7
$ws.Range("K77").Value = 3860923.5
$ws.Range("M77").Value = -3856555.5
$ws.Range("H132").Value = 4367.1577
$ws.Range("I132").Value = 3053.0967
$ws.Range("K132").Value = 9159.2901
$ws.Range("M132").Value = -6629.2901
$ws.Range("H136").Value = 6927.3335
$ws.Range("I136").Value = 4258.857
$ws.Range("K136").Value = 12776.571
$ws.Range("M136").Value = -10226.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 44444
$ws.Range("J62").Value = 44444
$ws.Range("L62").Value = 44444
$ws.Range("N62").Value = -45816
$ws.Range("H65").Value = 44444
$ws.Range("J65").Value = 44444
$ws.Range("L65").Value = 133332
$ws.Range("N65").Value = -140196
$ws.Range("H134").Value = 5766.2593
$ws.Range("I134").Value = 2409.125
$ws.Range("K134").Value = 7227.375
$ws.Range("M134").Value = -4692.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1735.1875
$ws.Range("J22").Value = 3613.3333
$ws.Range("L22").Value = 3613.3333
$ws.Range("N22").Value = -4313.3333
$ws.Range("H31").Value = 40004268
$ws.Range("I31").Value = 90909976
$ws.Range("K31").Value = 90909976
$ws.Range("M31").Value = -90909681
$ws.Range("H34").Value = 40004268
$ws.Range("I34").Value = 90909976
$ws.Range("K34").Value = 90909976
$ws.Range("M34").Value = -90909774
$ws.Range("H105").Value = 2658
$ws.Range("I105").Value = 2658
$ws.Range("K105").Value = 2658
$ws.Range("M105").Value = -911

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 111.65
$ws.Range("I33").Value = 109.75
$ws.Range("J33").Value = 119.25
$ws.Range("K33").Value = 658.5
$ws.Range("L33").Value = 715.5
$ws.Range("M33").Value = -375.5
$ws.Range("N33").Value = -1281.5
$ws.Range("H129").Value = 26316526
$ws.Range("J129").Value = 250000200
$ws.Range("L129").Value = 750000600
$ws.Range("N129").Value = -750010600

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 9921.799999999999
$ws.Range("I43").Value = 2220
$ws.Range("J43").Value = 17623.6
$ws.Range("K43").Value = 2220
$ws.Range("L43").Value = 17623.6
$ws.Range("M43").Value = -2069
$ws.Range("N43").Value = -17925.6
$ws.Range("H70").Value = 15714.549
$ws.Range("I70").Value = 12344.091
$ws.Range("K70").Value = 12344.091
$ws.Range("M70").Value = -12074.091
$ws.Range("H73").Value = 15714.549
$ws.Range("I73").Value = 12344.091
$ws.Range("K73").Value = 12344.091
$ws.Range("M73").Value = -11408.091
$ws.Range("H93").Value = 46997.75
$ws.Range("J93").Value = 46997.75
$ws.Range("L93").Value = 46997.75
$ws.Range("N93").Value = -50741.75
$ws.Range("H102").Value = 1724.1
$ws.Range("I102").Value = 1893.0625
$ws.Range("J102").Value = 1048.25
$ws.Range("K102").Value = 1893.0625
$ws.Range("L102").Value = 1048.25
$ws.Range("M102").Value = -271.0625
$ws.Range("N102").Value = -4292.25
$ws.Range("H122").Value = 8325.823
$ws.Range("I122").Value = 9538.571
$ws.Range("K122").Value = 28615.713
$ws.Range("M122").Value = -26165.713
$ws.Range("H126").Value = 2482.1333
$ws.Range("I126").Value = 2482.1333
$ws.Range("K126").Value = 7446.3999
$ws.Range("M126").Value = -4976.3999
$ws.Range("H132").Value = 5502.853
$ws.Range("I132").Value = 4691.516
$ws.Range("K132").Value = 14074.548
$ws.Range("M132").Value = -11544.548

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 168766.67
$ws.Range("I7").Value = 168766.67
$ws.Range("K7").Value = 168766.67
$ws.Range("M7").Value = -168654.67
$ws.Range("H61").Value = 1228.7142
$ws.Range("I61").Value = 1121.4445
$ws.Range("K61").Value = 1121.4445
$ws.Range("M61").Value = -919.4445000000001
$ws.Range("H113").Value = 1228.7142
$ws.Range("I113").Value = 1121.4445
$ws.Range("K113").Value = 1121.4445
$ws.Range("M113").Value = 1048.5555
$ws.Range("H126").Value = 168766.67
$ws.Range("I126").Value = 168766.67
$ws.Range("K126").Value = 506300.01
$ws.Range("M126").Value = -503830.01
$ws.Range("H132").Value = 9999.5
$ws.Range("I132").Value = 3499.5
$ws.Range("K132").Value = 10498.5
$ws.Range("M132").Value = -7968.5
$ws.Range("H136").Value = 6043.95
$ws.Range("I136").Value = 4100.4
$ws.Range("K136").Value = 12301.2
$ws.Range("M136").Value = -9751.199999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 72243.164
$ws.Range("J121").Value = 72243.164
$ws.Range("L121").Value = 72243.164
$ws.Range("N121").Value = -75737.164
$ws.Range("H126").Value = 4834.1763
$ws.Range("I126").Value = 3108.6365
$ws.Range("J126").Value = 7997.6665
$ws.Range("K126").Value = 9325.9095
$ws.Range("L126").Value = 23992.9995
$ws.Range("M126").Value = -6855.9095
$ws.Range("N126").Value = -28932.9995
$ws.Range("H132").Value = 6650.7393
$ws.Range("I132").Value = 3638.2727
$ws.Range("J132").Value = 9412.166999999999
$ws.Range("K132").Value = 10914.8181
$ws.Range("L132").Value = 28236.501
$ws.Range("M132").Value = -8384.8181
$ws.Range("N132").Value = -33296.501
